$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 4998.8
$ws.Range("I2").Value = 1666
$ws.Range("K2").Value = 1666
$ws.Range("M2").Value = -1553

$ws.Range("H40").Value = 2192.111
$ws.Range("J40").Value = 2191.125
$ws.Range("L40").Value = 2191.125
$ws.Range("N40").Value = -2541.125

$ws.Range("H52").Value = 1666.6666
$ws.Range("I52").Value = 1666.6666
$ws.Range("K52").Value = 4999.9998
$ws.Range("M52").Value = -4839.9998

$ws.Range("H112").Value = 1238.6774
$ws.Range("J112").Value = 1755.5
$ws.Range("L112").Value = 5266.5
$ws.Range("N112").Value = -7482.5

$ws.Range("H132").Value = 2578.3333
$ws.Range("I132").Value = 1796.6428
$ws.Range("K132").Value = 5389.928400000001
$ws.Range("M132").Value = -2859.928400000001

$ws.Range("H137").Value = 1776.4166
$ws.Range("I137").Value = 1620.7142
$ws.Range("K137").Value = 4862.142599999999
$ws.Range("M137").Value = -2312.142599999999

$ws.Range("H138").Value = 2684.0635
$ws.Range("J138").Value = 2950.2654
$ws.Range("L138").Value = 8850.796200000001
$ws.Range("N138").Value = -19130.7962

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 695.7826
$ws.Range("I2").Value = 1037.2307
$ws.Range("K2").Value = 1037.2307
$ws.Range("M2").Value = -924.2307000000001

$ws.Range("H32").Value = 10543
$ws.Range("I32").Value = 10543
$ws.Range("K32").Value = 10543
$ws.Range("M32").Value = -10256

$ws.Range("H45").Value = 1931.5
$ws.Range("I45").Value = 1998.8
$ws.Range("J45").Value = 1595
$ws.Range("K45").Value = 1998.8
$ws.Range("L45").Value = 1595
$ws.Range("M45").Value = -1621.8
$ws.Range("N45").Value = -2349

$ws.Range("H61").Value = 2330.5386
$ws.Range("I61").Value = 1829.0714
$ws.Range("J61").Value = 2915.5833
$ws.Range("K61").Value = 1829.0714
$ws.Range("L61").Value = 2915.5833
$ws.Range("M61").Value = -1617.0714
$ws.Range("N61").Value = -3339.5833

$ws.Range("H74").Value = 22216620
$ws.Range("I74").Value = 28562368
$ws.Range("K74").Value = 28562368
$ws.Range("M74").Value = -28561494

$ws.Range("H77").Value = 22216620
$ws.Range("I77").Value = 28562368
$ws.Range("K77").Value = 142811840
$ws.Range("M77").Value = -142807472

$ws.Range("H102").Value = 1349.5
$ws.Range("I102").Value = 1349.5
$ws.Range("K102").Value = 1349.5
$ws.Range("M102").Value = 272.5

$ws.Range("H112").Value = 50000
$ws.Range("J112").Value = 50000
$ws.Range("L112").Value = 50000
$ws.Range("N112").Value = -52954

$ws.Range("H116").Value = 695.7826
$ws.Range("I116").Value = 1037.2307
$ws.Range("K116").Value = 1037.2307
$ws.Range("M116").Value = 1256.7693

$ws.Range("H136").Value = 2330.5386
$ws.Range("I136").Value = 1829.0714
$ws.Range("J136").Value = 2915.5833
$ws.Range("K136").Value = 5487.2142
$ws.Range("L136").Value = 8746.749899999999
$ws.Range("M136").Value = -2937.2142
$ws.Range("N136").Value = -13846.7499

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 695.7826
$ws.Range("I3").Value = 1037.2307
$ws.Range("K3").Value = 1037.2307
$ws.Range("M3").Value = -923.2307000000001

$ws.Range("H31").Value = 5678
$ws.Range("J31").Value = 5678
$ws.Range("L31").Value = 5678
$ws.Range("N31").Value = -6182

$ws.Range("H99").Value = 884.75
$ws.Range("I99").Value = 885.3333
$ws.Range("J99").Value = 883
$ws.Range("K99").Value = 885.3333
$ws.Range("L99").Value = 883
$ws.Range("M99").Value = 612.6667
$ws.Range("N99").Value = -3879

$ws.Range("H134").Value = 1570.6571
$ws.Range("I134").Value = 1389.4517
$ws.Range("K134").Value = 4168.355100000001
$ws.Range("M134").Value = -1633.355100000001

$ws.Range("H135").Value = 174998
$ws.Range("J135").Value = 174998
$ws.Range("L135").Value = 174998
$ws.Range("N135").Value = -185138

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1325
$ws.Range("J16").Value = 996
$ws.Range("L16").Value = 996
$ws.Range("N16").Value = -1570

$ws.Range("H113").Value = 1325
$ws.Range("J113").Value = 996
$ws.Range("L113").Value = 996
$ws.Range("N113").Value = -5336

$ws.Range("H132").Value = 2525.0435
$ws.Range("I132").Value = 1967.6428
$ws.Range("J132").Value = 3392.111
$ws.Range("K132").Value = 5902.928400000001
$ws.Range("L132").Value = 10176.333
$ws.Range("M132").Value = -3372.928400000001
$ws.Range("N132").Value = -15236.333

$ws.Range("H134").Value = 2848.4285
$ws.Range("I134").Value = 2976.5334
$ws.Range("J134").Value = 2528.1667
$ws.Range("K134").Value = 8929.600199999999
$ws.Range("L134").Value = 7584.500100000001
$ws.Range("M134").Value = -6394.600199999999
$ws.Range("N134").Value = -12654.5001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 0
$ws.Range("J39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("N39").ClearContents()

$ws.Range("H40").Value = 350
$ws.Range("I40").Value = 275
$ws.Range("J40").Value = 500
$ws.Range("K40").Value = 1100
$ws.Range("L40").Value = 2000
$ws.Range("M40").Value = -1031
$ws.Range("N40").Value = -2138

$ws.Range("H41").Value = 2975
$ws.Range("I41").Value = 2975
$ws.Range("K41").Value = 8925
$ws.Range("M41").Value = -8587

$ws.Range("H104").Value = 2186.625
$ws.Range("I104").Value = 1200
$ws.Range("J104").Value = 2327.5715
$ws.Range("K104").Value = 3600
$ws.Range("L104").Value = 6982.7145
$ws.Range("M104").Value = -979
$ws.Range("N104").Value = -12224.7145

$ws.Range("H107").Value = 1761.1666
$ws.Range("J107").Value = 116.75
$ws.Range("L107").Value = 350.25
$ws.Range("N107").Value = -4190.25

$ws.Range("H109").Value = 2057.6
$ws.Range("I109").Value = 1822
$ws.Range("J109").Value = 3000
$ws.Range("K109").Value = 5466
$ws.Range("L109").Value = 9000
$ws.Range("M109").Value = -4426
$ws.Range("N109").Value = -11080

$ws.Range("H131").Value = 1221.7142
$ws.Range("J131").Value = 1832.3334
$ws.Range("L131").Value = 5497.0002
$ws.Range("N131").Value = -15577.0002

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 20840.666
$ws.Range("I43").Value = 10017
$ws.Range("J43").Value = 26252.5
$ws.Range("K43").Value = 10017
$ws.Range("L43").Value = 26252.5
$ws.Range("M43").Value = -9866
$ws.Range("N43").Value = -26554.5

$ws.Range("H46").Value = 34938.75
$ws.Range("J46").Value = 34938.75
$ws.Range("L46").Value = 34938.75
$ws.Range("N46").Value = -35250.75

$ws.Range("H80").Value = 4869.5
$ws.Range("I80").Value = 4370.875
$ws.Range("J80").Value = 5534.3335
$ws.Range("K80").Value = 4370.875
$ws.Range("L80").Value = 5534.3335
$ws.Range("M80").Value = -3372.875
$ws.Range("N80").Value = -7530.3335

$ws.Range("H83").Value = 4869.5
$ws.Range("I83").Value = 4370.875
$ws.Range("J83").Value = 5534.3335
$ws.Range("K83").Value = 21854.375
$ws.Range("L83").Value = 27671.6675
$ws.Range("M83").Value = -16862.375
$ws.Range("N83").Value = -37655.6675

$ws.Range("H122").Value = 2614.1667
$ws.Range("I122").Value = 1937
$ws.Range("K122").Value = 5811
$ws.Range("M122").Value = -3361

$ws.Range("H126").Value = 1262.9286
$ws.Range("I126").Value = 970.7778
$ws.Range("J126").Value = 1788.8
$ws.Range("K126").Value = 2912.3334
$ws.Range("L126").Value = 5366.4
$ws.Range("M126").Value = -442.3334
$ws.Range("N126").Value = -10306.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 450
$ws.Range("I16").Value = 450
$ws.Range("K16").Value = 450
$ws.Range("M16").Value = -280

$ws.Range("H46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("N46").ClearContents()

$ws.Range("H68").Value = 2955.5
$ws.Range("I68").Value = 2944.5
$ws.Range("J68").Value = 2988.5
$ws.Range("K68").Value = 2944.5
$ws.Range("L68").Value = 2988.5
$ws.Range("M68").Value = -2195.5
$ws.Range("N68").Value = -4486.5

$ws.Range("H71").Value = 2955.5
$ws.Range("I71").Value = 2944.5
$ws.Range("J71").Value = 2988.5
$ws.Range("K71").Value = 14722.5
$ws.Range("L71").Value = 14942.5
$ws.Range("M71").Value = -10978.5
$ws.Range("N71").Value = -22430.5

$ws.Range("H82").Value = 1857.2
$ws.Range("I82").Value = 1843.5
$ws.Range("J82").Value = 1866.3334
$ws.Range("K82").Value = 1843.5
$ws.Range("L82").Value = 1866.3334
$ws.Range("M82").Value = -1482.5
$ws.Range("N82").Value = -2588.3334

$ws.Range("H85").Value = 1857.2
$ws.Range("I85").Value = 1843.5
$ws.Range("J85").Value = 1866.3334
$ws.Range("K85").Value = 1843.5
$ws.Range("L85").Value = 1866.3334
$ws.Range("M85").Value = -595.5
$ws.Range("N85").Value = -4362.3334

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H80").Value = 32499
$ws.Range("J80").Value = 32499
$ws.Range("L80").Value = 32499
$ws.Range("N80").Value = -34495

$ws.Range("H83").Value = 32499
$ws.Range("J83").Value = 32499
$ws.Range("L83").Value = 97497
$ws.Range("N83").Value = -107481

$ws.Range("H132").Value = 3250.125
$ws.Range("I132").Value = 2865.7778
$ws.Range("J132").Value = 3744.2856
$ws.Range("K132").Value = 8597.3334
$ws.Range("L132").Value = 11232.8568
$ws.Range("M132").Value = -6067.3334
$ws.Range("N132").Value = -16292.8568
